$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column (Price) cells keep their literal text representation
# instead of being auto-converted to numbers (which would drop formatting
# such as trailing zeros, or reformat dotted "thousand" separators).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.234.19"
$ws.Range("E2").Value = "  +1.82%  "
$ws.Range("D3").Value = "1.796.81"
$ws.Range("E3").Value = "  +3.14%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "334.63"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "0.4508"
$ws.Range("E7").Value = "  +15.95%  "
$ws.Range("D8").Value = "0.3710"
$ws.Range("E8").Value = "  +10.23%  "
$ws.Range("D9").Value = "45.10"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").Value = "1.142"
$ws.Range("E10").Value = "  +4.01%  "
$ws.Range("D11").Value = "0.07572"
$ws.Range("E11").Value = "  +6.11%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "22.39"
$ws.Range("E13").Value = "  +3.18%  "
$ws.Range("D14").Value = "6.294"
$ws.Range("E14").Value = "  +3.93%  "
$ws.Range("D15").Value = "7.481"
$ws.Range("E15").Value = "  +7.95%  "
$ws.Range("D16").Value = "1.794.22"
$ws.Range("E16").Value = "  +3.03%  "
$ws.Range("D17").Value = "0.00001089"
$ws.Range("E17").Value = "  +4.22%  "
$ws.Range("D18").Value = "0.06743"
$ws.Range("E18").Value = "  +2.47%  "
$ws.Range("D19").Value = "81.13"
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "17.48"
$ws.Range("E21").Value = "  +4.77%  "
$ws.Range("D22").Value = "6.367"
$ws.Range("E22").Value = "  +3.66%  "
$ws.Range("D23").Value = "28.218.24"
$ws.Range("E23").Value = "  +1.70%  "
$ws.Range("D24").Value = "11.81"
$ws.Range("E24").Value = "  +2.97%  "
$ws.Range("D25").Value = "2.414"
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("D26").Value = "20.56"
$ws.Range("E26").Value = "  +4.57%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "2.365"
$ws.Range("E27").Value = "  +4.26%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "151.82"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("D29").Value = "1.997.81"
$ws.Range("E29").Value = "  +3.01%  "
$ws.Range("D30").Value = "132.84"
$ws.Range("E30").Value = "  +4.20%  "
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("D33").Value = "0.09407"
$ws.Range("E33").Value = "  +8.18%  "
$ws.Range("D34").Value = "5.802"
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("D35").Value = "0.2361"
$ws.Range("E35").Value = "  +13.79%  "
$ws.Range("D36").Value = "12.06"
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.02341"
$ws.Range("E37").Value = "  +3.91%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.06306"
$ws.Range("E38").Value = "  +4.13%  "
$ws.Range("D39").Value = "5.205"
$ws.Range("E39").Value = "  +2.80%  "
$ws.Range("D40").Value = "0.6565"
$ws.Range("E40").Value = "  +3.10%  "
$ws.Range("D41").Value = "8.366"
$ws.Range("E41").Value = "  +7.07%  "
$ws.Range("D42").Value = "1.483"
$ws.Range("E42").Value = "  -1.69%  "
$ws.Range("D43").Value = "1.213"
$ws.Range("E43").Value = "  +2.28%  "
$ws.Range("D44").Value = "14.26"
$ws.Range("E44").Value = "  +5.49%  "
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.6085"
$ws.Range("E46").Value = "  +3.32%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.827"
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("D48").Value = "130.17"
$ws.Range("E48").Value = "  +4.06%  "
$ws.Range("D49").Value = "2.033"
$ws.Range("E49").Value = "  +3.77%  "
$ws.Range("D50").Value = "0.07118"
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("D51").Value = "1.160"
$ws.Range("E51").Value = "  +1.92%  "
